$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.454.94"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "3.544.03"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "581.85"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "173.04"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.538.81"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").Value = "6.75"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").Value = "0.585"
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("D13").Value = "47.63"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("E14").Value = "  -4.67%  "
$ws.Range("D15").Value = "4.114.21"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").Value = "8.57"
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").Value = "631.23"
$ws.Range("E17").Value = "  -5.81%  "
$ws.Range("D18").Value = "69.512.12"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "3.538.11"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").Value = "0.124"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "0.892"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.00"
$ws.Range("E24").Value = "  -6.52%  "
$ws.Range("D25").Value = "97.57"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "5.84"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "2.65"
$ws.Range("E29").Value = "  -5.04%  "
$ws.Range("D30").Value = "9.36"
$ws.Range("E30").Value = "  -6.33%  "
$ws.Range("D31").Value = "32.96"
$ws.Range("E31").Value = "  -5.35%  "
$ws.Range("E32").Value = "  -5.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.60"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("D34").Value = "1.35"
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").Value = "7.02"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("D36").Value = "638.98"
$ws.Range("E36").Value = "  +9.37%  "
$ws.Range("D37").Value = "10.82"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").Value = "  -11.37%  "
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("D40").Value = "57.45"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "0.0458"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "0.137"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").Value = "3.402.09"
$ws.Range("E44").Value = "  -5.68%  "
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("D46").Value = "0.0₃0706"
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("D47").Value = "32.87"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("E48").Value = "  -5.17%  "
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "132.99"
$ws.Range("E51").Value = "  -1.57%  "
